# Updated symbol list on Sun Dec 18 06:26:42 UTC 2022 with GitHub Actions
#
# This script applies a refreshed snapshot of crypto prices (column D),
# together with a few label/link updates (columns B, C, E) to Sheet1 of
# the workbook, matching the upstream data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of silently re-parsing numeric-looking strings as numbers
    # (all data cells in this sheet - including "Price" - are text cells).
    $r = $ws.Range($Cell)
    $r.Value = "'" + $Text
    # The quote-prefix entry bumps the cell to a new style; restore the
    # original (default) style so only the value itself changes.
    $r.Style = "Normal"
}

# ---- Price (column D) refresh ----
Set-TextValue "D2"  "246.65"
Set-TextValue "D3"  "22.79"
Set-TextValue "D4"  "5.610"
Set-TextValue "D6"  "3.405"
Set-TextValue "D7"  "6.475"
Set-TextValue "D8"  "0.8026"
Set-TextValue "D9"  "1.066"
Set-TextValue "D11" "0.07488"
Set-TextValue "D12" "0.03188"
Set-TextValue "D13" "0.02972"
Set-TextValue "D14" "0.09262"
Set-TextValue "D15" "0.001665"
Set-TextValue "D17" "0.04693"
Set-TextValue "D18" "0.0005745"
Set-TextValue "D19" "0.006265"
Set-TextValue "D20" "0.001056"
Set-TextValue "D21" "0.003813"
Set-TextValue "D23" "0.0004604"
Set-TextValue "D25" "2.121"
Set-TextValue "D27" "0.1278"
Set-TextValue "D40" "0.04186"

# ---- Row 18: "Worst in 24h" badge moves off of One (ONE) ----
$ws.Range("E18").Value = "17OneONE"

# ---- Rows 41-43: KickToken / BKEXToken / CEJI reshuffle position+rank ----
# Row 41 becomes BKEXToken (was KickToken)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1048"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 becomes CEJI (was BKEXToken)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002972"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes KickToken (was CEJI), and now carries the "Worst in 24h" badge
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003251"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# ---- Remaining price (column D) refresh ----
Set-TextValue "D44" "0.009803"
Set-TextValue "D45" "0.00005634"
Set-TextValue "D47" "0.6806"
Set-TextValue "D48" "0.02883"
